# Add data for 2021-11-26
# Update the "through" date references and bump/insert the latest-month
# (column B, header row 1) counts plus a handful of the November 2020
# (column M) counts that were also revised, for the neighborhoods touched
# by this data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab to reflect the new "through" date.
$ws.Name = "Through 2021-11-18"

# Update the header label for the current/latest month column (column B).
$ws.Range("B1").Value = "November 2021 (through November 18)"

# --- Column B (November 2021, through Nov 18) updates ---
$ws.Range("B3").Value = 5     # Garfield Park
$ws.Range("B6").Value = 7     # West Town
$ws.Range("B44").Value = 1    # East Village (new)
$ws.Range("B59").Value = 3    # Little Village
$ws.Range("B88").Value = 1    # Old Town (new)

# --- Column M (November 2020) updates ---
$ws.Range("M9").Value = 4     # Auburn Gresham
$ws.Range("M10").Value = 3    # Lower West Side
$ws.Range("M12").Value = 5    # Grand Boulevard
$ws.Range("M15").Value = 2    # Calumet Heights
$ws.Range("M16").Value = 4    # Washington Heights
$ws.Range("M18").Value = 1    # Loop (new)
$ws.Range("M34").Value = 1    # Hyde Park (new)
$ws.Range("M46").Value = 1    # Fuller Park (new)
$ws.Range("M96").Value = 2    # Ukrainian Village

# --- Other scattered single-cell updates ---
$ws.Range("BP4").Value = 3    # Austin, November 2015
$ws.Range("AT5").Value = 2    # Humboldt Park, November 2017
$ws.Range("X6").Value = 2     # West Town, November 2019
$ws.Range("AT8").Value = 4    # Englewood, November 2017
$ws.Range("AI19").Value = 1   # Chinatown, November 2018 (new)
$ws.Range("AI25").Value = 4   # Uptown, November 2018
